# daily auto push: 2026-01-13 09:38 UTC
# Insert a new data row for 2026/01/13 (time slot 17) into the daily log
# table. The new row is inserted immediately before the existing row 628
# (the block that starts the 2026/12/29 data), pushing rows 628:669 down
# to 629:670.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 628, shifting existing rows 628:669 down
# to 629:670 (dimension grows from A1:D669 to A1:D670).
$ws.Rows.Item(628).Insert()

# Populate the newly inserted row 628 with the new data point.
# Columns A and B hold text (date / weekday-in-Japanese) stored as plain
# strings in the original file, so force text storage on A628 (which
# looks like a date and would otherwise be auto-converted to a date
# serial number), then drop the formatting again so the cell ends up
# with the same "no explicit style" look as its neighboring data rows.
$ws.Cells.Item(628, 1).NumberFormat = "@"
$ws.Cells.Item(628, 1).Value = "2026/01/13"
$ws.Cells.Item(628, 1).ClearFormats()

$ws.Cells.Item(628, 2).Value = "火"
$ws.Cells.Item(628, 3).Value = 17
$ws.Cells.Item(628, 4).Value = 201
